# Natural_Gas_Production.xlsx — add EF 2020 data points + TC Energy (Coalbed Methane etc.)
# ratings for the "pq" (Power Query result) sheet.
#
# The "pq" sheet backs an ExternalData/query table ("natural_gas_prod") with
# columns Year | Production Type | Production (BCf/d). Historically it held
# 2015-2019 rows for each of 5 production types (Solution, Non Associated,
# Tight, Shale, Coalbed Methane). This edit appends a 2020 row after each
# production type's 2015-2019 block (5 new rows total), then grows the
# query table / named range / dimension to cover the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pq")
$tbl = $ws.ListObjects.Item(1)

# --- Insert the five new "2020" rows, each right after its production
#     type's existing 2015-2019 block, shifting everything below it down. ---

# Solution 2020 -> new row 7 (after Solution 2015-2019 in rows 2-6)
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value2 = 2020
$ws.Range("B7").Value2 = "Solution"
$ws.Range("C7").Value2 = 2.15

# Non Associated 2020 -> new row 13 (after Non Associated 2015-2019, now rows 8-12)
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Value2 = 2020
$ws.Range("B13").Value2 = "Non Associated"
$ws.Range("C13").Value2 = 2.19

# Tight 2020 -> new row 19 (after Tight 2015-2019, now rows 14-18)
$ws.Rows.Item(19).Insert()
$ws.Range("A19").Value2 = 2020
$ws.Range("B19").Value2 = "Tight"
$ws.Range("C19").Value2 = 10.42

# Shale 2020 -> new row 25 (after Shale 2015-2019, now rows 20-24)
$ws.Rows.Item(25).Insert()
$ws.Range("A25").Value2 = 2020
$ws.Range("B25").Value2 = "Shale"
$ws.Range("C25").Value2 = 0.51

# Coalbed Methane 2020 -> new row 31 (after Coalbed Methane 2015-2019, now rows 26-30)
$ws.Rows.Item(31).Insert()
$ws.Range("A31").Value2 = 2020
$ws.Range("B31").Value2 = "Coalbed Methane"
$ws.Range("C31").Value2 = 0.46

# --- Grow the query table / autofilter to the new extent A1:C31. ---
$tbl.Resize($ws.Range("A1:C31"))

# --- Keep the workbook-level ExternalData_1 hidden name in sync with the
#     query table's new extent. ---
$wb.Names.Item("ExternalData_1").RefersTo = "=pq!`$A`$1:`$C`$31"

# --- Match the author's final cursor position on the pq sheet. ---
$ws.Activate()
$ws.Range("E6").Select()
